$d = $word.ActiveDocument

$replacements = @(
    @("2025-02-25 Tuesday", "2025-02-26 Wednesday"),
    @("837×8=", "149×8="),
    @("238×7=", "358×3="),
    @("343×5=", "653×2="),
    @("459×5=", "411×2="),
    @("516×9=", "344×7="),
    @("776×8=", "622×9="),
    @("398×9=", "115×3="),
    @("468×2=", "289×6="),
    @("538×5=", "162×4="),
    @("192×8=", "931×3="),
    @("707×7=", "297×6="),
    @("269×9=", "331×7="),
    @("472×3=", "674×4="),
    @("881×9=", "854×6="),
    @("151×5=", "542×3="),
    @("569×4=", "174×8="),
    @("174×2=", "154×4="),
    @("195×5=", "655×3="),
    @("519×3=", "911×4="),
    @("872×4=", "831×9="),
    @("193×4=", "788×9="),
    @("965×3=", "314×3="),
    @("358×4=", "794×3="),
    @("966×5=", "792×4="),
    @("889×3=", "629×3=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
